# Merging of VWAP data and brand data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Linde row (row 22) with Industry / Sector / Ticker
$ws.Range("E22").Value = "Speciality Chemicals"
$ws.Range("F22").Value = "Basic Materials"
$ws.Range("G22").Value = "LIN.F"

# That row now wraps onto two lines like the other populated rows
$ws.Rows.Item(22).RowHeight = 28

# Normalize the Volkswagen ticker text (drop the stray non-breaking space / mixed run formatting)
$ws.Range("G30").Value = "VOW3.F"
$ws.Range("G31").Value = "VOW3.F"

# Restore the view to where the edit was made
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E22").Select()
